# feat: add 2022-Q3 data
#
# 1) "总计" (overview) sheet: insert a new 2022-Q3 summary row above the
#    existing 2021-Q2 summary row (which shifts down to row 3).
# 2) Add a brand-new "2022-Q3" worksheet (placed right after "总计",
#    before the existing "2021-Q2" sheet) holding the per-fund detail
#    rows for the new quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# -----------------------------------------------------------------
# 1) "总计" sheet updates
# -----------------------------------------------------------------

# Give row 3's index cell (A3) the same "index column" style as A2
# before filling it in, so it keeps the bold/centered/bordered look.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Push the existing 2021-Q2 totals down into row 3.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

# Write the new 2022-Q3 totals into row 2.
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

# -----------------------------------------------------------------
# 2) New "2022-Q3" worksheet (inserted right after "总计")
# -----------------------------------------------------------------

$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# --- header row (row 1), matching the bold/bordered header style ---
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$total.Range("B1:D1").Copy()
$q3.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("B1:D1").Copy()
$q3.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("B1").Copy()
$q3.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 2: fund 096001 ---
$total.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$q3.Range("A2").Value = 0

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "096001"
$q3.Range("B2").ClearFormats()

$q3.Range("C2").NumberFormat = "@"
$q3.Range("C2").Value = "大成标普500等权重指数（QDII）人民币"
$q3.Range("C2").ClearFormats()

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "3.08"
$q3.Range("D2").ClearFormats()

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "93.16"
$q3.Range("E2").ClearFormats()

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "0.23"
$q3.Range("F2").ClearFormats()

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0071"
$q3.Range("G2").ClearFormats()

$q3.Range("H2").Value = 3

# --- row 3: fund 013404 ---
$total.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$q3.Range("A3").Value = 1

$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "013404"
$q3.Range("B3").ClearFormats()

$q3.Range("C3").NumberFormat = "@"
$q3.Range("C3").Value = "大成标普500等权重指数（QDII）美元"
$q3.Range("C3").ClearFormats()

$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "3.08"
$q3.Range("D3").ClearFormats()

$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "93.16"
$q3.Range("E3").ClearFormats()

$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "0.23"
$q3.Range("F3").ClearFormats()

$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0071"
$q3.Range("G3").ClearFormats()

$q3.Range("H3").Value = 3

# Keep "总计" as the active/selected sheet, as in the original file.
$total.Activate()
